$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -1
    4  = 1
    5  = 3
    6  = 1
    7  = 2
    8  = -2
    9  = -7
    10 = 3
    11 = 0
    13 = 7
    14 = -3
    15 = -2
    16 = -1
    17 = -1
    18 = -1
    19 = -1
    20 = 1
    21 = 0
    22 = -1
    24 = 2
    25 = -1
    26 = 5
    28 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
